$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value2 = "ECs"
$ws.Cells.Item(2, 2).Value2 = "Ccl25"
$ws.Cells.Item(2, 3).Value2 = "Ackr4"
$ws.Cells.Item(2, 4).Value2 = "ECs"
$ws.Cells.Item(2, 5).Value2 = 3
$ws.Cells.Item(2, 6).Value2 = 1
$ws.Cells.Item(2, 7).Value2 = 5.934604333333334
$ws.Cells.Item(2, 8).Value2 = 17.803813
$ws.Cells.Item(2, 9).Value2 = 0.3081877218757661
$ws.Cells.Item(2, 10).Value2 = 0.3081877218757661
$ws.Cells.Item(2, 11).Value2 = 1
$ws.Cells.Item(2, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(2, 13).Value2 = 0.1143813333333333
$ws.Cells.Item(2, 14).Value2 = 0.343144
$ws.Cells.Item(2, 15).Value2 = 0.03414257747883775
$ws.Cells.Item(2, 16).Value2 = 0.03414257747883775
$ws.Cells.Item(2, 17).Value2 = 0.6788079564524445
$ws.Cells.Item(2, 18).Value2 = 6.109271608072
$ws.Cells.Item(2, 19).Value2 = 0.01052232317216984
$ws.Cells.Item(2, 20).Value2 = 0.01052232317216984

# Row 3
$ws.Cells.Item(3, 1).Value2 = "ECs"
$ws.Cells.Item(3, 2).Value2 = "Ccl25"
$ws.Cells.Item(3, 3).Value2 = "Ackr4"
$ws.Cells.Item(3, 4).Value2 = "FAPs"
$ws.Cells.Item(3, 5).Value2 = 3
$ws.Cells.Item(3, 6).Value2 = 1
$ws.Cells.Item(3, 7).Value2 = 5.934604333333334
$ws.Cells.Item(3, 8).Value2 = 17.803813
$ws.Cells.Item(3, 9).Value2 = 0.3081877218757661
$ws.Cells.Item(3, 10).Value2 = 0.3081877218757661
$ws.Cells.Item(3, 11).Value2 = 3
$ws.Cells.Item(3, 12).Value2 = 1
$ws.Cells.Item(3, 13).Value2 = 3.175982333333334
$ws.Cells.Item(3, 14).Value2 = 9.527947000000001
$ws.Cells.Item(3, 15).Value2 = 0.9480237703755849
$ws.Cells.Item(3, 16).Value2 = 0.9480237703755849
$ws.Cells.Item(3, 17).Value2 = 18.84819851799011
$ws.Cells.Item(3, 18).Value2 = 169.633786661911
$ws.Cells.Item(3, 19).Value2 = 0.2921692860761259
$ws.Cells.Item(3, 20).Value2 = 0.2921692860761259

# Row 4
$ws.Cells.Item(4, 1).Value2 = "ECs"
$ws.Cells.Item(4, 2).Value2 = "Ccl25"
$ws.Cells.Item(4, 3).Value2 = "Ackr4"
$ws.Cells.Item(4, 4).Value2 = "sCs"
$ws.Cells.Item(4, 5).Value2 = 3
$ws.Cells.Item(4, 6).Value2 = 1
$ws.Cells.Item(4, 7).Value2 = 5.934604333333334
$ws.Cells.Item(4, 8).Value2 = 17.803813
$ws.Cells.Item(4, 9).Value2 = 0.3081877218757661
$ws.Cells.Item(4, 10).Value2 = 0.3081877218757661
$ws.Cells.Item(4, 11).Value2 = 2
$ws.Cells.Item(4, 12).Value2 = 0.6666666666666666
$ws.Cells.Item(4, 13).Value2 = 0.05974466666666667
$ws.Cells.Item(4, 14).Value2 = 0.179234
$ws.Cells.Item(4, 15).Value2 = 0.01783365214557738
$ws.Cells.Item(4, 16).Value2 = 0.01783365214557738
$ws.Cells.Item(4, 17).Value2 = 0.3545609576935556
$ws.Cells.Item(4, 18).Value2 = 3.191048619242
$ws.Cells.Item(4, 19).Value2 = 0.005496112627470361
$ws.Cells.Item(4, 20).Value2 = 0.005496112627470361

# Row 5
$ws.Cells.Item(5, 1).Value2 = "FAPs"
$ws.Cells.Item(5, 2).Value2 = "Ccl25"
$ws.Cells.Item(5, 3).Value2 = "Ackr4"
$ws.Cells.Item(5, 4).Value2 = "ECs"
$ws.Cells.Item(5, 5).Value2 = 3
$ws.Cells.Item(5, 6).Value2 = 1
$ws.Cells.Item(5, 7).Value2 = 7.720664
$ws.Cells.Item(5, 8).Value2 = 23.161992
$ws.Cells.Item(5, 9).Value2 = 0.4009389195777736
$ws.Cells.Item(5, 10).Value2 = 0.4009389195777736
$ws.Cells.Item(5, 11).Value2 = 1
$ws.Cells.Item(5, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(5, 13).Value2 = 0.1143813333333333
$ws.Cells.Item(5, 14).Value2 = 0.343144
$ws.Cells.Item(5, 15).Value2 = 0.03414257747883775
$ws.Cells.Item(5, 16).Value2 = 0.03414257747883775
$ws.Cells.Item(5, 17).Value2 = 0.8830998425386667
$ws.Cells.Item(5, 18).Value2 = 7.947898582848
$ws.Cells.Item(5, 19).Value2 = 0.01368908812596563
$ws.Cells.Item(5, 20).Value2 = 0.01368908812596563

# Row 6
$ws.Cells.Item(6, 1).Value2 = "FAPs"
$ws.Cells.Item(6, 2).Value2 = "Ccl25"
$ws.Cells.Item(6, 3).Value2 = "Ackr4"
$ws.Cells.Item(6, 4).Value2 = "FAPs"
$ws.Cells.Item(6, 5).Value2 = 3
$ws.Cells.Item(6, 6).Value2 = 1
$ws.Cells.Item(6, 7).Value2 = 7.720664
$ws.Cells.Item(6, 8).Value2 = 23.161992
$ws.Cells.Item(6, 9).Value2 = 0.4009389195777736
$ws.Cells.Item(6, 10).Value2 = 0.4009389195777736
$ws.Cells.Item(6, 11).Value2 = 3
$ws.Cells.Item(6, 12).Value2 = 1
$ws.Cells.Item(6, 13).Value2 = 3.175982333333334
$ws.Cells.Item(6, 14).Value2 = 9.527947000000001
$ws.Cells.Item(6, 15).Value2 = 0.9480237703755849
$ws.Cells.Item(6, 16).Value2 = 0.9480237703755849
$ws.Cells.Item(6, 17).Value2 = 24.52069246560267
$ws.Cells.Item(6, 18).Value2 = 220.6862321904241
$ws.Cells.Item(6, 19).Value2 = 0.3800996262284343
$ws.Cells.Item(6, 20).Value2 = 0.3800996262284343

# Row 7
$ws.Cells.Item(7, 1).Value2 = "FAPs"
$ws.Cells.Item(7, 2).Value2 = "Ccl25"
$ws.Cells.Item(7, 3).Value2 = "Ackr4"
$ws.Cells.Item(7, 4).Value2 = "sCs"
$ws.Cells.Item(7, 5).Value2 = 3
$ws.Cells.Item(7, 6).Value2 = 1
$ws.Cells.Item(7, 7).Value2 = 7.720664
$ws.Cells.Item(7, 8).Value2 = 23.161992
$ws.Cells.Item(7, 9).Value2 = 0.4009389195777736
$ws.Cells.Item(7, 10).Value2 = 0.4009389195777736
$ws.Cells.Item(7, 11).Value2 = 2
$ws.Cells.Item(7, 12).Value2 = 0.6666666666666666
$ws.Cells.Item(7, 13).Value2 = 0.05974466666666667
$ws.Cells.Item(7, 14).Value2 = 0.179234
$ws.Cells.Item(7, 15).Value2 = 0.01783365214557738
$ws.Cells.Item(7, 16).Value2 = 0.01783365214557738
$ws.Cells.Item(7, 17).Value2 = 0.4612684971253334
$ws.Cells.Item(7, 18).Value2 = 4.151416474128
$ws.Cells.Item(7, 19).Value2 = 0.007150205223373639
$ws.Cells.Item(7, 20).Value2 = 0.007150205223373639

# Row 8
$ws.Cells.Item(8, 1).Value2 = "sCs"
$ws.Cells.Item(8, 2).Value2 = "Ccl25"
$ws.Cells.Item(8, 3).Value2 = "Ackr4"
$ws.Cells.Item(8, 4).Value2 = "ECs"
$ws.Cells.Item(8, 5).Value2 = 3
$ws.Cells.Item(8, 6).Value2 = 1
$ws.Cells.Item(8, 7).Value2 = 5.601191
$ws.Cells.Item(8, 8).Value2 = 16.803573
$ws.Cells.Item(8, 9).Value2 = 0.2908733585464604
$ws.Cells.Item(8, 10).Value2 = 0.2908733585464603
$ws.Cells.Item(8, 11).Value2 = 1
$ws.Cells.Item(8, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(8, 13).Value2 = 0.1143813333333333
$ws.Cells.Item(8, 14).Value2 = 0.343144
$ws.Cells.Item(8, 15).Value2 = 0.03414257747883775
$ws.Cells.Item(8, 16).Value2 = 0.03414257747883775
$ws.Cells.Item(8, 17).Value2 = 0.6406716948346667
$ws.Cells.Item(8, 18).Value2 = 5.766045253512
$ws.Cells.Item(8, 19).Value2 = 0.009931166180702276
$ws.Cells.Item(8, 20).Value2 = 0.009931166180702274

# Row 9
$ws.Cells.Item(9, 1).Value2 = "sCs"
$ws.Cells.Item(9, 2).Value2 = "Ccl25"
$ws.Cells.Item(9, 3).Value2 = "Ackr4"
$ws.Cells.Item(9, 4).Value2 = "FAPs"
$ws.Cells.Item(9, 5).Value2 = 3
$ws.Cells.Item(9, 6).Value2 = 1
$ws.Cells.Item(9, 7).Value2 = 5.601191
$ws.Cells.Item(9, 8).Value2 = 16.803573
$ws.Cells.Item(9, 9).Value2 = 0.2908733585464604
$ws.Cells.Item(9, 10).Value2 = 0.2908733585464603
$ws.Cells.Item(9, 11).Value2 = 3
$ws.Cells.Item(9, 12).Value2 = 1
$ws.Cells.Item(9, 13).Value2 = 3.175982333333334
$ws.Cells.Item(9, 14).Value2 = 9.527947000000001
$ws.Cells.Item(9, 15).Value2 = 0.9480237703755849
$ws.Cells.Item(9, 16).Value2 = 0.9480237703755849
$ws.Cells.Item(9, 17).Value2 = 17.78928366162567
$ws.Cells.Item(9, 18).Value2 = 160.103552954631
$ws.Cells.Item(9, 19).Value2 = 0.2757548580710248
$ws.Cells.Item(9, 20).Value2 = 0.2757548580710247

# Row 10
$ws.Cells.Item(10, 1).Value2 = "sCs"
$ws.Cells.Item(10, 2).Value2 = "Ccl25"
$ws.Cells.Item(10, 3).Value2 = "Ackr4"
$ws.Cells.Item(10, 4).Value2 = "sCs"
$ws.Cells.Item(10, 5).Value2 = 3
$ws.Cells.Item(10, 6).Value2 = 1
$ws.Cells.Item(10, 7).Value2 = 5.601191
$ws.Cells.Item(10, 8).Value2 = 16.803573
$ws.Cells.Item(10, 9).Value2 = 0.2908733585464604
$ws.Cells.Item(10, 10).Value2 = 0.2908733585464603
$ws.Cells.Item(10, 11).Value2 = 2
$ws.Cells.Item(10, 12).Value2 = 0.6666666666666666
$ws.Cells.Item(10, 13).Value2 = 0.05974466666666667
$ws.Cells.Item(10, 14).Value2 = 0.179234
$ws.Cells.Item(10, 15).Value2 = 0.01783365214557738
$ws.Cells.Item(10, 16).Value2 = 0.01783365214557738
$ws.Cells.Item(10, 17).Value2 = 0.3346412892313333
$ws.Cells.Item(10, 18).Value2 = 3.011771603082
$ws.Cells.Item(10, 19).Value2 = 0.005187334294733383
$ws.Cells.Item(10, 20).Value2 = 0.005187334294733382
